$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "Datos actualizados..." timestamp cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 2 de Octubre de 2020 a las 19:26"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 7508662
$ws.Range("C4").Value = 13991
$ws.Range("D4").Value = 4750427
$ws.Range("E4").Value = 2545300
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 275
$ws.Range("H4").Value = 212935

# Row 21 - Turquia
$ws.Range("B21").Value = 321512
$ws.Range("C21").Value = 1442
$ws.Range("D21").Value = 282657
$ws.Range("E21").Value = 30530
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 63
$ws.Range("H21").Value = 8325

# Row 56 - Chequia
$ws.Range("B56").Value = 76017
$ws.Range("C56").Value = 1762
$ws.Range("D56").Value = 34871
$ws.Range("E56").Value = 40447
$ws.Range("F56").Value = 0
$ws.Range("G56").Value = 21
$ws.Range("H56").Value = 699

# Row 68 was Paraguay, now becomes Libano (reordered alphabetically ahead of Paraguay)
$ws.Range("A68").Value = "Libano"
$ws.Range("B68").Value = 42173
$ws.Range("C68").Value = 1291
$ws.Range("D68").Value = 18379
$ws.Range("E68").Value = 23408
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 12
$ws.Range("H68").Value = 386

# Row 69 was Libano, now becomes Paraguay
$ws.Range("A69").Value = "Paraguay"
$ws.Range("B69").Value = 41799
$ws.Range("C69").Value = 0
$ws.Range("D69").Value = 25167
$ws.Range("E69").Value = 15763
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = 869
